# Insert a new data row at row 384 (this shifts existing rows 384..416 down to 385..417)
# and populate the newly inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).Insert()

$newRow = 384

$ws.Cells.Item($newRow, 1).Value = 10
$ws.Cells.Item($newRow, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value = 45132
$ws.Cells.Item($newRow, 5).Value = 9
$ws.Cells.Item($newRow, 6).Value = 100112052
$ws.Cells.Item($newRow, 7).Value = "Albahaca"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 45
$ws.Cells.Item($newRow, 11).Value = 6000
$ws.Cells.Item($newRow, 12).Value = 6000
$ws.Cells.Item($newRow, 13).Value = 6000
$ws.Cells.Item($newRow, 14).Value = "$/paquete"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 6000
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
